# Workbook "wb1" gains a sibling sheet "wb2" with a small date-formatting
# demo table, plus a new row of numbers on "wb1" itself. The new sheet
# becomes the active one.

$wb = $excel.ActiveWorkbook
$wb1 = $wb.Worksheets.Item(1)

# --- wb1: add a row of five numbers below the existing data -----------
$wb1.Range("B16").Value = 1
$wb1.Range("C16").Value = 2
$wb1.Range("D16").Value = 3
$wb1.Range("E16").Value = 4
$wb1.Range("F16").Value = 5

# --- add the new sheet right after wb1 ---------------------------------
$wb2 = $wb.Worksheets.Add($null, $wb1)
$wb2.Name = "wb2"

# Row 2: a real date serial formatted MM/DD/YY, with an explanatory note
$wb2.Range("B2").Value = 43586
$wb2.Range("B2").NumberFormat = "MM/DD/YY"
$wb2.Range("D2").Value = "// For whatever reason the first date is in MM/DD/YY"

# Row 3: a date typed in as literal (DD/MM/YYYY) text, same number format
# applied to the cell even though the content is text, not a serial
$wb2.Range("B3").Value = "27/12/2016"
$wb2.Range("B3").NumberFormat = "MM/DD/YY"
$wb2.Range("D3").WrapText = $true

# Row 4: another literal date-as-text
$wb2.Range("B4").Value = "23/07/1976"
$wb2.Range("B4").NumberFormat = "MM/DD/YY"

# Row 5: a second real date serial, with its own note
$wb2.Range("B5").Value = 31595
$wb2.Range("B5").NumberFormat = "MM/DD/YY"
$wb2.Range("D5").Value = "// So is this"

# Row 6: one more literal date-as-text, with a final note
$wb2.Range("B6").Value = "19/01/2038"
$wb2.Range("B6").NumberFormat = "MM/DD/YY"
$wb2.Range("D6").Value = "// This is a string"

# Row 7: blank cell, but carries the date style like its neighbours above
$wb2.Range("B7").NumberFormat = "MM/DD/YY"

# The new sheet ends up active, with B7 selected; wb1 keeps a stale
# selection at C19 from before the insert.
$null = $wb1.Range("C19").Select()
$null = $wb2.Range("B7").Select()
$null = $wb2.Activate()
